$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(6).Insert()
Write-Host "done insert"
for ($r = 1; $r -le 12; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    Write-Host ("Row $r : A=$a B=$b")
}
